$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = 0.8191621144995483
$ws.Range("C4").Value = 0.5990300658092822
$ws.Range("C5").Value = -0.05884812859031883
$ws.Range("C6").Value = -0.2477909811474888
$ws.Range("C7").Value = 0.08798946409716052
$ws.Range("C8").Value = 0.1903335571888334
$ws.Range("C9").Value = 0.08731620729673588
$ws.Range("C10").Value = -0.0001454401703124805
$ws.Range("C12").Value = 1.386943624147526
